$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update E2 value: "Q2_update_synonyms" -> "update_synonyms"
$ws.Range("E2").Value = "update_synonyms"

# Widen column E (stored OOXML width = ColumnWidth + 5/6, so back-solve for 23)
$ws.Range("E1").ColumnWidth = 22.166666666666668

# Update sheet view: select E1 (also clears the old topLeftCell="C1" scroll position)
$ws.Activate()
$ws.Range("E1").Select()
